$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 344.25
$ws.Range("J5").Value = 545
$ws.Range("L5").Value = 545
$ws.Range("N5").Value = -775
$ws.Range("H11").Value = 274.75
$ws.Range("I11").Value = 274.75
$ws.Range("K11").Value = 274.75
$ws.Range("M11").Value = -134.75
$ws.Range("H19").Value = 950.44446
$ws.Range("J19").Value = 969.625
$ws.Range("L19").Value = 969.625
$ws.Range("N19").Value = -1319.625
$ws.Range("H41").Value = 262.1111
$ws.Range("I41").Value = 232.375
$ws.Range("K41").Value = 232.375
$ws.Range("M41").Value = 207.625
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H98").Value = 1733.3334
$ws.Range("I98").Value = 1100
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 1100
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 398
$ws.Range("N98").Value = -5996
$ws.Range("H106").Value = 4496.6665
$ws.Range("I106").Value = 4496.6665
$ws.Range("K106").Value = 4496.6665
$ws.Range("M106").Value = -3865.6665
$ws.Range("H122").Value = 1733.3334
$ws.Range("I122").Value = 1100
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3300
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -850
$ws.Range("N122").Value = -13900
$ws.Range("H129").Value = 1099.5
$ws.Range("I129").Value = 1099.5
$ws.Range("K129").Value = 3298.5
$ws.Range("M129").Value = 1701.5
$ws.Range("H137").Value = 924.8333
$ws.Range("I137").Value = 869.8
$ws.Range("K137").Value = 2609.4
$ws.Range("M137").Value = -59.39999999999964
$ws.Range("H138").Value = 6500
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2085.4443
$ws.Range("I61").Value = 2085.4443
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2085.4443
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1873.4443
$ws.Range("N61").ClearContents()
$ws.Range("I74").Value = 2999
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2999
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2125
$ws.Range("N74").ClearContents()
$ws.Range("I77").Value = 2999
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 14995
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -10627
$ws.Range("N77").ClearContents()
$ws.Range("H97").Value = 1299.2963
$ws.Range("I97").Value = 1071.7
$ws.Range("J97").Value = 1949.5714
$ws.Range("K97").Value = 1071.7
$ws.Range("L97").Value = 1949.5714
$ws.Range("M97").Value = -575.7
$ws.Range("N97").Value = -2941.5714
$ws.Range("H136").Value = 2085.4443
$ws.Range("I136").Value = 2085.4443
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6256.3329
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3706.3329
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3134.9
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H134").Value = 3749
$ws.Range("I134").Value = 3749
$ws.Range("K134").Value = 11247
$ws.Range("M134").Value = -8712

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1899.3334
$ws.Range("I31").Value = 1749
$ws.Range("J31").Value = 2200
$ws.Range("K31").Value = 1749
$ws.Range("L31").Value = 2200
$ws.Range("M31").Value = -1454
$ws.Range("N31").Value = -2790
$ws.Range("H32").Value = 3237.7144
$ws.Range("I32").Value = 3237.7144
$ws.Range("K32").Value = 3237.7144
$ws.Range("M32").Value = -2921.7144
$ws.Range("H34").Value = 1899.3334
$ws.Range("I34").Value = 1749
$ws.Range("J34").Value = 2200
$ws.Range("K34").Value = 1749
$ws.Range("L34").Value = 2200
$ws.Range("M34").Value = -1547
$ws.Range("N34").Value = -2604
$ws.Range("H107").Value = 3005.9167
$ws.Range("I107").Value = 2825
$ws.Range("J107").Value = 3096.375
$ws.Range("K107").Value = 2825
$ws.Range("L107").Value = 3096.375
$ws.Range("M107").Value = -905
$ws.Range("N107").Value = -6936.375
$ws.Range("H134").Value = 2264.353
$ws.Range("I134").Value = 2264.353
$ws.Range("K134").Value = 6793.059
$ws.Range("M134").Value = -4258.059

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1964.6666
$ws.Range("I86").Value = 2225
$ws.Range("J86").Value = 1444
$ws.Range("K86").Value = 6675
$ws.Range("L86").Value = 4332
$ws.Range("M86").Value = -5489
$ws.Range("N86").Value = -6704
$ws.Range("H89").Value = 1964.6666
$ws.Range("I89").Value = 2225
$ws.Range("J89").Value = 1444
$ws.Range("K89").Value = 20025
$ws.Range("L89").Value = 12996
$ws.Range("M89").Value = -14097
$ws.Range("N89").Value = -24852
$ws.Range("H98").Value = 296.5
$ws.Range("I98").Value = 296.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 889.5
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 608.5
$ws.Range("N98").ClearContents()
$ws.Range("H104").Value = 2029
$ws.Range("J104").Value = 2029
$ws.Range("L104").Value = 6087
$ws.Range("N104").Value = -11329
$ws.Range("H109").Value = 1000
$ws.Range("I109").Value = 1000
$ws.Range("K109").Value = 3000
$ws.Range("M109").Value = -1960
$ws.Range("H115").Value = 900
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 900
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 2700
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -5050
$ws.Range("H121").Value = 580.1667
$ws.Range("I121").Value = 289.6
$ws.Range("J121").Value = 2033
$ws.Range("K121").Value = 868.8000000000001
$ws.Range("L121").Value = 6099
$ws.Range("M121").Value = 441.1999999999999
$ws.Range("N121").Value = -8719
$ws.Range("H139").Value = 3391.25
$ws.Range("I139").Value = 3391.25
$ws.Range("K139").Value = 10173.75
$ws.Range("M139").Value = -5033.75
$ws.Range("H140").Value = 649.6667
$ws.Range("I140").Value = 649.6667
$ws.Range("K140").Value = 1949.0001
$ws.Range("M140").Value = 3230.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 162.83333
$ws.Range("I2").Value = 84.166664
$ws.Range("J2").Value = 241.5
$ws.Range("K2").Value = 84.166664
$ws.Range("L2").Value = 241.5
$ws.Range("M2").Value = 28.833336
$ws.Range("N2").Value = -467.5
$ws.Range("H102").Value = 2400.2
$ws.Range("I102").Value = 2125.25
$ws.Range("K102").Value = 2125.25
$ws.Range("M102").Value = -503.25
$ws.Range("H107").Value = 3654.125
$ws.Range("I107").Value = 633.0909
$ws.Range("K107").Value = 633.0909
$ws.Range("M107").Value = 1286.9091
$ws.Range("H122").Value = 1139.4
$ws.Range("I122").Value = 1099.6666
$ws.Range("J122").Value = 1199
$ws.Range("K122").Value = 3298.9998
$ws.Range("L122").Value = 3597
$ws.Range("M122").Value = -848.9998000000001
$ws.Range("N122").Value = -8497
$ws.Range("H126").Value = 3490.7856
$ws.Range("J126").Value = 4226.5
$ws.Range("L126").Value = 12679.5
$ws.Range("N126").Value = -17619.5
$ws.Range("H132").Value = 2293.4167
$ws.Range("I132").Value = 2242.2222
$ws.Range("K132").Value = 6726.6666
$ws.Range("M132").Value = -4196.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1995
$ws.Range("I32").Value = 1995
$ws.Range("K32").Value = 1995
$ws.Range("M32").Value = -1678
$ws.Range("H46").Value = 4472.04
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 5180.1
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 5180.1
$ws.Range("M46").Value = -3812
$ws.Range("N46").Value = -5556.1
$ws.Range("H93").Value = 685.3333
$ws.Range("J93").Value = 766.3333
$ws.Range("L93").Value = 766.3333
$ws.Range("N93").Value = -3262.3333
$ws.Range("H100").Value = 34200
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H132").Value = 5219.375
$ws.Range("I132").Value = 4626.5
$ws.Range("J132").Value = 6998
$ws.Range("K132").Value = 13879.5
$ws.Range("L132").Value = 20994
$ws.Range("M132").Value = -11349.5
$ws.Range("N132").Value = -26054

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 69950
$ws.Range("J46").Value = 69950
$ws.Range("L46").Value = 69950
$ws.Range("N46").Value = -70412
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H96").Value = 1789.9375
$ws.Range("I96").Value = 1798.1818
$ws.Range("J96").Value = 1771.8
$ws.Range("K96").Value = 1798.1818
$ws.Range("L96").Value = 1771.8
$ws.Range("M96").Value = -425.1818000000001
$ws.Range("N96").Value = -4517.8
$ws.Range("H100").Value = 3320535.2
$ws.Range("J100").Value = 2707.8333
$ws.Range("L100").Value = 5415.6666
$ws.Range("N100").Value = -6497.6666
$ws.Range("H109").Value = 22000
$ws.Range("J109").Value = 22000
$ws.Range("L109").Value = 22000
$ws.Range("N109").Value = -24774
$ws.Range("H134").Value = 69950
$ws.Range("J134").Value = 69950
$ws.Range("L134").Value = 209850
$ws.Range("N134").Value = -214920
